# Update countries & provincias Spain
# Refreshes the COVID-19 dashboard data: updates the "last updated" timestamp,
# updates case counts for several countries, and re-establishes the
# column-B (Casos totales) descending sort order for the rows whose counts
# changed enough to move past a neighboring row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 12:17"

# --- Banglades (row 18) ---
$ws.Range("B18").Value = 332970
$ws.Range("C18").Value = 1892
$ws.Range("D18").Value = 233550
$ws.Range("E18").Value = 94786
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = 4634

# --- Malasia (row 97) ---
$ws.Range("B97").Value = 9628
$ws.Range("C97").Value = 45
$ws.Range("D97").Value = 9167
$ws.Range("E97").Value = 333

# --- Finlandia (row 103) ---
$ws.Range("B103").Value = 8469
$ws.Range("C103").Value = 39
$ws.Range("E103").Value = 632

# --- Rows 126-131: Jamaica, Mayotte, Somalia, Siria, Eslovenia, Gambia ---
# Eslovenia's refreshed count (3389) now outranks Mayotte/Somalia/Siria, so it
# moves up to row 127 (right after Jamaica); Mayotte/Somalia/Siria each shift
# down one row, keeping their own (unchanged) figures.
$ws.Range("A127").Value = "Eslovenia"
$ws.Range("B127").Value = 3389
$ws.Range("C127").Value = 77
$ws.Range("D127").Value = 2620
$ws.Range("E127").Value = 634
$ws.Range("H127").Value = 135

$ws.Range("A128").Value = "Mayotte"
$ws.Range("B128").Value = 3374
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 2964
$ws.Range("E128").Value = 370
$ws.Range("H128").Value = 40

$ws.Range("A129").Value = "Somalia"
$ws.Range("B129").Value = 3371
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 2738
$ws.Range("E129").Value = 536
$ws.Range("H129").Value = 97

$ws.Range("A130").Value = "Siria"
$ws.Range("B130").Value = 3351
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 780
$ws.Range("E130").Value = 2428
$ws.Range("H130").Value = 143

# --- Sri Lanka (row 133) ---
$ws.Range("D133").Value = 2955
$ws.Range("E133").Value = 180

# --- Rows 214-215: Montserrat / Islas Malvinas swap (tied total cases) ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
